# The commit updates the table on slide 16 to use a different built-in
# PowerPoint table style (swapping the table style applied to the
# income/expenditure summary table).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)

if ($shp.HasTable) {
    $tbl = $shp.Table
    $tbl.ApplyStyle("{7765F754-B564-4037-8E7B-6B17DF8D3FAE}")
}
